$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Cover-letter date: "November 14, 2022" -> "December 12, 2022"
#    (typed as separate edits, so it lands in Word as 4 runs:
#     "December" | " 1" | "2" | ", 2022")
# ---------------------------------------------------------------
$d.Content.Find.Execute("November 14, 2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "December 12, 2022", 2) | Out-Null

$dateRng = $d.Content
$dateRng.Find.Execute("December 12, 2022") | Out-Null
$dateStart = $dateRng.Start

$p1 = $d.Range($dateStart, $dateStart + 8)          # "December"
$p1.Bold = $true
$p1.Bold = $false

$p2 = $d.Range($dateStart + 8, $dateStart + 10)     # " 1"
$p2.Bold = $true
$p2.Bold = $false

$p3 = $d.Range($dateStart + 10, $dateStart + 11)    # "2"
$p3.Bold = $true
$p3.Bold = $false

# ---------------------------------------------------------------
# 2) Manuscript title: "Lorem ipsum" -> full title text, still
#    bold+italic, split into 3 runs with spell-check markers
#    around "Riffomonas".
# ---------------------------------------------------------------
$d.Content.Find.Execute("Lorem ipsum", $true, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "The Riffomonas YouTube Channel: An Educational Resource to Foster Reproducible Research Practices", `
                         2) | Out-Null

$titleRng = $d.Content
$titleRng.Find.Execute("The Riffomonas YouTube Channel: An Educational Resource to Foster Reproducible Research Practices") | Out-Null
$titleStart = $titleRng.Start

$t1 = $d.Range($titleStart, $titleStart + 4)        # "The "
$t1.Bold = $true
$t1.Bold = $false

$t2 = $d.Range($titleStart + 4, $titleStart + 14)   # "Riffomonas"
$t2.Bold = $true
$t2.Bold = $false

# ---------------------------------------------------------------
# 3) Subscriber count sentence: "...currently has 11,000
#    subscribers..." -> "...currently has more than 11,300
#    subscribers...", split into 5 runs.
# ---------------------------------------------------------------
$oldSentence = "The resource that this manuscript announces is a YouTube channel that teaches researchers how to engage in reproducible research practices. The channel currently has 11,000 subscribers and 285 videos. I am frequently asked by viewers whether there is something that they can cite to give credit for the skills they have developed on the channel that have become components of their research."
$newSentence = "The resource that this manuscript announces is a YouTube channel that teaches researchers how to engage in reproducible research practices. The channel currently has more than 11,300 subscribers and 285 videos. I am frequently asked by viewers whether there is something that they can cite to give credit for the skills they have developed on the channel that have become components of their research."

$d.Content.Find.Execute($oldSentence, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newSentence, 2) | Out-Null

$subRng = $d.Content
$subRng.Find.Execute($newSentence) | Out-Null
$subStart = $subRng.Start

$s1 = $d.Range($subStart, $subStart + 166)          # "...currently has "
$s1.Bold = $true
$s1.Bold = $false

$s2 = $d.Range($subStart + 166, $subStart + 176)    # "more than "
$s2.Bold = $true
$s2.Bold = $false

$s3 = $d.Range($subStart + 176, $subStart + 179)    # "11,"
$s3.Bold = $true
$s3.Bold = $false

$s4 = $d.Range($subStart + 179, $subStart + 180)    # "3"
$s4.Bold = $true
$s4.Bold = $false

# ---------------------------------------------------------------
# 4) Paragraph alignment: explicitly set "Align Left" on every
#    paragraph from the salutation through to the end of the
#    letter (they currently inherit "justify" from the Normal
#    style).
# ---------------------------------------------------------------
$greetRng = $d.Content
$greetRng.Find.Execute("Dear Dr. Newton,") | Out-Null
$bodyStart = $greetRng.Start

$endRng = $d.Content
$endRng.Find.Execute("Professor") | Out-Null
$bodyEnd = $endRng.End

$body = $d.Range($bodyStart, $bodyEnd)
foreach ($para in $body.Paragraphs) {
    $para.Range.ParagraphFormat.Alignment = 0
}

# ---------------------------------------------------------------
# 5) styles.xml: DefaultParagraphFont is no longer semi-hidden.
# ---------------------------------------------------------------
$styles = $d.Styles
$dpf = $styles.Item("Default Paragraph Font")
$dpf.Hidden = $false
